# Weekly data update: append Wk47 rows (928-954) to "Weekly Expenditure" sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly expenditure rows for Wk47 (columns A-H)
$newRows = @(
    @("Wk47", "SAP", 11155143, "PTS-1136 300-001519-015 Semiconductor", "4000pcs", "Fishes", "17/11/2025", 7440),
    @("Wk47", "SAP", 11155143, "PTS-1136 300-001519-015 Semiconductor", "200pcs", "Fishes", "18/11/2025", 372),
    @("Wk47", "SAP", "XS-PTS-0876", "HX 2067 HPN CRC CO CONTACT CLEANER", "4pcs", "Fishes", "18/11/2025", 90.92),
    @("Wk47", "SAP", 11156537, "PTS-1150 300-01698-010 Semiconductor", "200pcs", "Fishes", "18/11/2025", 468),
    @("Wk47", "SAP", "XS-PTS-0469", "70508.201 10# BULB MFG: PHILIPS", "1pcs", "Fishes", "18/11/2025", 131.88),
    @("Wk47", "SAP", "XS-PTS-0098", "CDUJB6-10D AIR CYLINDER", "8pcs", "Fishes", "18/11/2025", 186.31),
    @("Wk47", "SAP", "XS-PTS-0356", "GP-570D84A-03 Hyperspace Semicon Probes", "1000pcs", "Sihl", "19/11/2025", 2100),
    @("Wk47", "SAP", "XS-PTS-0124", "70700.545 COMRRESSION SPRING (Input &", "2pcs", "Fishes", "19/11/2025", 17.18),
    @("Wk47", "SAP", "XS-PTS-1026", "10618.381 SECONDARY BOARD X1767", "2pcs", "Fishes", "19/11/2025", 2311.74),
    @("Wk47", "SAP", 11157203, "PTS-1160 10416.080 2 metre 0.635 MM HIGH", "2pcs", "Fishes", "20/11/2025", 1081.12),
    @("Wk47", "SAP", 11155949, "PTS-1146 70902.631 X2544 ALIGNER & NEST", "4pcs", "Fishes", "20/11/2025", 2889.48),
    @("Wk47", "SAP", 11156537, "PTS-1150 300-01698-010 Semiconductor", "200pcs", "Fishes", "20/11/2025", 468),
    @("Wk47", "SAP", "XS-PTS-1025", "10618.38 PRIMARY PCB BOARD ASSY X1767", "4pcs", "Fishes", "20/11/2025", 421),
    @("Wk47", "SAP", "XS-PTS-1021", "10819.01410 METER HIGH FLEX HYBRID", "1pcs", "Fishes", "20/11/2025", 361.39),
    @("Wk47", "Expense", "Expense", "PVC Wire 2C 2.5mm Green / Yellow 100meter", "2rolls", "Fishes", "21/11/2025", 78.01000000000005),
    @("Wk47", "Expense", "Expense", "Cable Lug (Spade) 2.5mm (100pcs Per Pack)", "2packs", "Fishes", "21/11/2025", 11.32),
    @("Wk47", "Expense", "Expense", "Snap-Off Blade (Retractable)", "2pcs", "Fishes", "21/11/2025", 26),
    @("Wk47", "Expense", "Expense", "SSF Signal Cable (2M)", "21pcs", "Fishes", "21/11/2025", 1009.1),
    @("Wk47", "SAP", "XS-PTS-0952", "10819.396 Yamaha Cable_Input & Output", "2pcs", "Fishes", "21/11/2025", 3275.36),
    @("Wk47", "SAP", "XS-PTS-1026", "10618.381 SECONDARY BOARD X1767", "3pcs", "Fishes", "21/11/2025", 3467.61),
    @("Wk47", "SAP", "XS-PTS-0953", "10819.292 8 METER HIGH FLEX HYBRID", "1pcs", "Fishes", "21/11/2025", 413.18),
    @("Wk47", "SAP", 11154686, "PTS-1129 TW.50.15.98.JV.016.01 X2544", "3pcs", "Fishes", "21/11/2025", 1583.35),
    @("Wk47", "SAP", 11155143, "PTS-1136 300-001519-015 Semiconductor", "1500pcs", "Fishes", "21/11/2025", 2790),
    @("Wk47", "Expense", "Expense", "Lencent Adaptor GaN III", "2pcs", "Lisa", "21/11/2025", 90.22),
    @("Wk47", "Expense", "Expense", "6'' Inch Sciss", "2pcs", "Lisa", "21/11/2025", 8.56),
    @("Wk47", "SAP", "XS-PTS-0299", "10018.004 RENISHAW LINEAR ENCODER", "1pcs", "Fishes", "22/11/2025", 4088.98),
    @("Wk47", "SAP", 11151248, "PTS-1071 TW.50.15.FI.0S.151.00 X2637", "3pcs", "Fishes", "22/11/2025", 1666.68)
)

$startRow = 928
$templateRow = 927
$lastRow = $startRow + $newRows.Count - 1

# Pre-format the new rows by inheriting styles/number formats from the last existing data row
$ws.Range("A$templateRow`:H$templateRow").Copy() | Out-Null
$ws.Range("A$startRow`:H$lastRow").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Fill in the row details (columns B-H) first, matching the order the data was entered
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    $ws.Range("B$r").Value = $rowData[1]
    $ws.Range("C$r").Value = $rowData[2]
    $ws.Range("D$r").Value = $rowData[3]
    $ws.Range("E$r").Value = $rowData[4]
    $ws.Range("F$r").Value = $rowData[5]
    $ws.Range("G$r").Value = $rowData[6]
    $ws.Range("H$r").Value = $rowData[7]
}

# Then fill the WorkWeek column (A) last, as it was filled down across the whole block
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    $ws.Range("A$r").Value = $rowData[0]
}

# Extend the AutoFilter range to cover the newly added rows
$ws.AutoFilterMode = $false
$ws.Range("A1:H$lastRow").AutoFilter() | Out-Null

# Keep the _FilterDatabase defined name in sync with the new filter range
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "='Weekly Expenditure'!`$A`$1:`$H`$$lastRow"
    }
}

# Move the active selection (matches the author's last cursor position)
$ws.Range("D17").Select() | Out-Null
